$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9-38: Koppen climate classification labels + hex colors
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = 'Af Tropical, rainforest'
$ws.Range("C9").Value = '0000FF'
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = 'Am Tropical, monsoon'
$ws.Range("C10").Value = '0078FF'
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = 'Aw Tropical, savannah'
$ws.Range("C11").Value = '46AAF'
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = 'BWh Arid, desert, hot'
$ws.Range("C12").Value = 'FF0000'
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = 'BWk Arid, desert, cold'
$ws.Range("C13").Value = 'FF9696'
$ws.Range("A14").Value = 6
$ws.Range("B14").Value = 'BSh Arid, steppe, hot'
$ws.Range("C14").Value = 'F5A500'
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = 'BSk Arid, steppe, cold'
$ws.Range("C15").Value = 'FFDC64'
$ws.Range("A16").Value = 8
$ws.Range("B16").Value = 'Csa Temperate, dry summer, hot summer'
$ws.Range("C16").Value = 'FFFF00'
$ws.Range("A17").Value = 9
$ws.Range("B17").Value = 'Csb Temperate, dry summer, warm summer'
$ws.Range("C17").Value = 'C8C800'
$ws.Range("A18").Value = 10
$ws.Range("B18").Value = 'Csc Temperate, dry summer, cold summer'
$ws.Range("C18").Value = 969600
$ws.Range("A19").Value = 11
$ws.Range("B19").Value = 'Cwa Temperate, dry winter, hot summer'
$ws.Range("C19").Value = '96FF96'
$ws.Range("A20").Value = 12
$ws.Range("B20").Value = 'Cwb Temperate, dry winter, warm summer'
$ws.Range("C20").Value = '64C864'
$ws.Range("A21").Value = 13
$ws.Range("B21").Value = 'Cwc Temperate, dry winter, cold summer'
$ws.Range("C21").Value = 329632
$ws.Range("A22").Value = 14
$ws.Range("B22").Value = 'Cfa Temperate, no dry season, hot summer'
$ws.Range("C22").Value = 'C8FF50'
$ws.Range("A23").Value = 15
$ws.Range("B23").Value = 'Cfb Temperate, no dry season, warm summer'
$ws.Range("C23").Value = '64FF50'
$ws.Range("A24").Value = 16
$ws.Range("B24").Value = 'Cfc Temperate, no dry season, cold summer'
$ws.Range("C24").Value = '32C800'
$ws.Range("A25").Value = 17
$ws.Range("B25").Value = 'Dsa Cold, dry summer, hot summer'
$ws.Range("C25").Value = 'FF00FF'
$ws.Range("A26").Value = 18
$ws.Range("B26").Value = 'Dsb Cold, dry summer, warm summer'
$ws.Range("C26").Value = 'C800C8'
$ws.Range("A27").Value = 19
$ws.Range("B27").Value = 'Dsc Cold, dry summer, cold summer'
$ws.Range("C27").Value = 963296
$ws.Range("A28").Value = 20
$ws.Range("B28").Value = 'Dsd Cold, dry summer, very cold winter'
$ws.Range("C28").Value = 966496
$ws.Range("A29").Value = 21
$ws.Range("B29").Value = 'Dwa Cold, dry winter, hot summer'
$ws.Range("C29").Value = 'AAAF'
$ws.Range("A30").Value = 22
$ws.Range("B30").Value = 'Dwb Cold, dry winter, warm summer'
$ws.Range("C30").Value = '5A78DC'
$ws.Range("A31").Value = 23
$ws.Range("B31").Value = 'Dwc Cold, dry winter, cold summer'
$ws.Range("C31").Value = '4B50B4'
$ws.Range("A32").Value = 24
$ws.Range("B32").Value = 'Dwd Cold, dry winter, very cold winter'
$ws.Range("C32").Value = 320087
$ws.Range("A33").Value = 25
$ws.Range("B33").Value = 'Dfa Cold, no dry season, hot summer'
$ws.Range("C33").Value = '00FFFF'
$ws.Range("A34").Value = 26
$ws.Range("B34").Value = 'Dfb Cold, no dry season, warm summer'
$ws.Range("C34").Value = '37C8FF'
$ws.Range("A35").Value = 27
$ws.Range("B35").Value = 'Dfc Cold, no dry season, cold summer'
$ws.Range("C35").Value = '007D7D'
$ws.Range("A36").Value = 28
$ws.Range("B36").Value = 'Dfd Cold, no dry season, very cold winter'
$ws.Range("C36").Value = '00465F'
$ws.Range("A37").Value = 29
$ws.Range("B37").Value = 'ET Polar, tundra'
$ws.Range("C37").Value = 'B2B2B2'
$ws.Range("A38").Value = 30
$ws.Range("B38").Value = 'EF Polar, frost'
$ws.Range("C38").Value = 666666

# Extend the CONCATENATE-based markdown-swatch formula down through the new rows
$ws.Range("D9:D38").Formula = '=CONCATENATE("![#",C9,"](https://placehold.co/15x15/",C9,"/",C9,".png)")'

# Column B needs to be widened to fit the longer climate-zone descriptions
$ws.Columns("B").ColumnWidth = 38

# Leave the selection where the editing session left it
$ws.Range("G16").Select()
